$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New data rows (13-19)
$data = @(
    @("2024-06-10_R_e.dat", "R", 0, 45453, 4,  1.8, 20, 202),
    @("2024-06-12_W_e.dat", "W", 0, 45455, 5,  1.8, 20, 202),
    @("2024-06-13_I_e.dat", "I", 0, 45456, 5,  1.8, 20, 202),
    @("2024-06-13_L_e.dat", "L", 0, 45456, 10, 1.8, 10, 202),
    @("2024-06-14_D_e.dat", "D", 0, 45457, 10, 1.8, 10, 202),
    @("2024-06-17_I_e.dat", "I", 0, 45460, 10, 1.8, 10, 202),
    @("2024-06-18_D_e.dat", "D", 0, 45461, 8,  1.8, 10, 202)
)

$r = 13
foreach ($row in $data) {
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    $ws.Cells.Item($r, 5).Value = $row[4]
    $ws.Cells.Item($r, 6).Value = $row[5]
    $ws.Cells.Item($r, 7).Value = $row[6]
    $ws.Cells.Item($r, 8).Value = $row[7]
    $r = $r + 1
}

# Copy the date format from D2 onto the new date cells so it reuses the
# existing style (numFmtId 14) instead of minting a new custom format.
$ws.Range("D2").Copy()
$ws.Range("D13:D19").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Column A width change (target stored width 22.1640625; the host quantizes
# ColumnWidth to 1/6-character steps, so 21.33 is the input that lands on
# the closest reachable stored width, 22.1666...)
$ws.Columns.Item(1).ColumnWidth = 21.33

# Selection moves to G19
$ws.Range("G19").Select() | Out-Null
